$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for columns I and J, matching the style of the
# existing header row (bold font, thin border, centered/top aligned).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108
$ws.Range("I1:J1").VerticalAlignment = -4160
$ws.Range("I1:J1").Borders.LineStyle = 1

# Per-row values for new columns I (I0) and J (IF).
$values = @(
    @{Row=2; I=7; J=8},
    @{Row=3; I=1; J=4},
    @{Row=4; I=1; J=5},
    @{Row=5; I=1; J=6},
    @{Row=6; I=1; J=2},
    @{Row=7; I=1; J=4},
    @{Row=8; I=1; J=5},
    @{Row=9; I=1; J=7},
    @{Row=10; I=1; J=6},
    @{Row=11; I=1; J=4},
    @{Row=12; I=1; J=5},
    @{Row=13; I=1; J=6},
    @{Row=14; I=1; J=6},
    @{Row=15; I=1; J=6},
    @{Row=16; I=1; J=6},
    @{Row=17; I=1; J=6},
    @{Row=18; I=1; J=6},
    @{Row=19; I=1; J=6},
    @{Row=20; I=1; J=6},
    @{Row=21; I=1; J=5},
    @{Row=22; I=1; J=5},
    @{Row=23; I=1; J=6},
    @{Row=24; I=1; J=7},
    @{Row=25; I=1; J=6},
    @{Row=26; I=1; J=5},
    @{Row=27; I=1; J=6},
    @{Row=28; I=1; J=6},
    @{Row=29; I=1; J=5},
    @{Row=30; I=1; J=7},
    @{Row=31; I=1; J=5},
    @{Row=32; I=4; J=6},
    @{Row=33; I=1; J=2}
)

foreach ($entry in $values) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I
    $ws.Cells.Item($r, 10).Value = $entry.J
}
